$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text or $null if unchanged, new Volume(1h) (E) percentage text
$updates = @(
    @(2, "97.603.89", "+3.75%"),
    @(3, "3.339.23", "+8.48%"),
    @(4, $null, "+0.05%"),
    @(5, "257.60", "+9.93%"),
    @(6, "619.95", "+1.83%"),
    @(7, $null, "+2.24%"),
    @(8, "0.386", "+2.29%"),
    @(9, $null, "+0.07%"),
    @(10, "3.338.25", "+8.54%"),
    @(11, "0.795", "-3.16%"),
    @(12, "0.200", "+1.76%"),
    @(13, "97.237.26", "+3.56%"),
    @(14, "35.55", "+4.63%"),
    @(15, "0.0000246", "+2.89%"),
    @(16, "3.956.94", "+8.52%"),
    @(17, "5.53", "+4.88%"),
    @(18, "3.333.22", "+8.17%"),
    @(19, "3.58", "-1.91%"),
    @(20, "14.98", "+3.01%"),
    @(21, "482.74", "+9.38%"),
    @(22, $null, "+9.18%"),
    @(23, "5.82", "+1.58%"),
    @(24, "9.23", "+4.68%"),
    @(25, "5.64", "+1.88%"),
    @(26, "88.03", "+3.60%"),
    @(27, "12.11", "+1.30%"),
    @(28, "3.522.06", "+8.31%"),
    @(29, $null, "+0.04%"),
    @(30, $null, "+3.54%"),
    @(31, "0.239", "-3.05%"),
    @(32, $null, "-0.71%"),
    @(33, $null, "+0.66%"),
    @(34, "9.25", "+0.91%"),
    @(35, "27.40", "+6.98%"),
    @(36, "7.42", "-4.10%"),
    @(37, $null, "-4.24%"),
    @(38, "511.40", "+9.51%"),
    @(39, $null, "+3.33%"),
    @(40, "24.82", "+3.47%"),
    @(41, "0.450", "+1.24%"),
    @(42, $null, "+1.11%"),
    @(43, "3.33", "+6.78%"),
    @(44, "3.54", "-4.79%"),
    @(45, "0.792", "+17.14%"),
    @(46, $null, "+0.03%"),
    @(47, "160.96", "+0.84%"),
    @(48, "1.92", "+3.50%"),
    @(49, $null, "+7.23%"),
    @(50, "45.55", "+4.26%"),
    @(51, "4.51", "+4.98%"),
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceText = $u[1]
    $volumeText = $u[2]

    if ($null -ne $priceText) {
        # Force plain text so numeric-looking prices ("257.60", "0.200", ...) keep
        # their exact original formatting instead of being auto-coerced to a number.
        $priceCell = $ws.Range("D$row")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceText
        $priceCell.Style = "Normal"
    }

    $ws.Range("E$row").Value = "  $volumeText  "
}
